$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 13, shifting existing rows 13..89 down to 14..90
$ws.Rows.Item(13).Insert()

# Fill in the newly inserted row 13
$ws.Range("A13").Value = "dct:creator"
$ws.Range("B13").Value = "Jitka"

# Update B12 with the new creator identifier
$ws.Range("B12").Value = "https://orcid.org/0000-0002-0454-4289"
